$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cells for the season record columns.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold/bordered header style used by the rest of row 1 (e.g. A1).
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the season record (Wins/Losses/Ties) for every data row.
for ($row = 2; $row -le 46; $row++) {
    $ws.Cells.Item($row, 30).Value = 96
    $ws.Cells.Item($row, 31).Value = 66
    $ws.Cells.Item($row, 32).Value = 0
}
